{"js": "const replacements = [\n  [\"947\u00f72=473, 1\", \"246\u00f79=27, 3\"],\n  [\"768\u00f78=96, 0\", \"165\u00f72=82, 1\"],\n  [\"331\u00f77=47, 2\", \"523\u00f73=174, 1\"],\n  [\"927\u00f77=132, 3\", \"321\u00f74=80, 1\"],\n  [\"352\u00f75=70, 2\", \"134\u00f72=67, 0\"],\n  [\"354\u00f75=70, 4\", \"582\u00f72=291, 0\"],\n  [\"557\u00f78=69, 5\", \"885\u00f76=147, 3\"],\n  [\"581\u00f75=116, 1\", \"596\u00f79=66, 2\"],\n  [\"350\u00f74=87, 2\", \"331\u00f77=47, 2\"],\n  [\"614\u00f73=204, 2\", \"783\u00f77=111, 6\"],\n  [\"589\u00f79=65, 4\", \"214\u00f77=30, 4\"],\n  [\"162\u00f76=27, 0\", \"841\u00f76=140, 1\"],\n  [\"104\u00f79=11, 5\", \"690\u00f79=76, 6\"],\n  [\"632\u00f77=90, 2\", \"227\u00f78=28, 3\"],\n  [\"546\u00f75=109, 1\", \"733\u00f72=366, 1\"],\n  [\"760\u00f76=126, 4\", \"696\u00f73=232, 0\"],\n  [\"647\u00f78=80, 7\", \"505\u00f79=56, 1\"],\n  [\"855\u00f78=106, 7\", \"263\u00f72=131, 1\"],\n  [\"369\u00f77=52, 5\", \"829\u00f76=138, 1\"],\n  [\"137\u00f76=22, 5\", \"790\u00f72=395, 0\"],\n  [\"781\u00f72=390, 1\", \"559\u00f75=111, 4\"],\n  [\"970\u00f72=485, 0\", \"253\u00f78=31, 5\"],\n  [\"344\u00f77=49, 1\", \"576\u00f76=96, 0\"],\n  [\"379\u00f77=54, 1\", \"408\u00f74=102, 0\"],\n  [\"800\u00f79=88, 8\", \"245\u00f78=30, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @('947\u00f72=473, 1', '246\u00f79=27, 3'),\n    @('768\u00f78=96, 0', '165\u00f72=82, 1'),\n    @('331\u00f77=47, 2', '523\u00f73=174, 1'),\n    @('927\u00f77=132, 3', '321\u00f74=80, 1'),\n    @('352\u00f75=70, 2', '134\u00f72=67, 0'),\n    @('354\u00f75=70, 4', '582\u00f72=291, 0'),\n    @('557\u00f78=69, 5', '885\u00f76=147, 3'),\n    @('581\u00f75=116, 1', '596\u00f79=66, 2'),\n    @('350\u00f74=87, 2', '331\u00f77=47, 2'),\n    @('614\u00f73=204, 2', '783\u00f77=111, 6'),\n    @('589\u00f79=65, 4', '214\u00f77=30, 4'),\n    @('162\u00f76=27, 0', '841\u00f76=140, 1'),\n    @('104\u00f79=11, 5', '690\u00f79=76, 6'),\n    @('632\u00f77=90, 2', '227\u00f78=28, 3'),\n    @('546\u00f75=109, 1', '733\u00f72=366, 1'),\n    @('760\u00f76=126, 4', '696\u00f73=232, 0'),\n    @('647\u00f78=80, 7', '505\u00f79=56, 1'),\n    @('855\u00f78=106, 7', '263\u00f72=131, 1'),\n    @('369\u00f77=52, 5', '829\u00f76=138, 1'),\n    @('137\u00f76=22, 5', '790\u00f72=395, 0'),\n    @('781\u00f72=390, 1', '559\u00f75=111, 4'),\n    @('970\u00f72=485, 0', '253\u00f78=31, 5'),\n    @('344\u00f77=49, 1', '576\u00f76=96, 0'),\n    @('379\u00f77=54, 1', '408\u00f74=102, 0'),\n    @('800\u00f79=88, 8', '245\u00f78=30, 5'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
